# Refresh cryptos list with latest market data (GitHub Actions scheduled run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "26.214.91"
$ws.Range("E2").Value = "  -6.22%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "1.667.56"
$ws.Range("E3").Value = "  -4.18%  "

# Row 5 (BNB)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.89%  "

# Row 6 (XRP)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5061"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -12.55%  "

# Row 7 (USDC)
$ws.Range("E7").Value = "  +0.43%  "

# Row 8 (Cardano)
$ws.Range("E8").Value = "  -3.31%  "

# Row 9 (Dogecoin)
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06353"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.16%  "

# Row 10 (Solana)
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.17%  "

# Row 11 (TRON)
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07356"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.66%  "

# Row 12 (WrappedEther -> Polkadot)
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.539"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.58%  "

# Row 13 (Polkadot -> WrappedEther)
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.663.28"
$ws.Range("E13").Value = "  -4.53%  "

# Row 14 (Polygon)
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5802"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.76%  "

# Row 15 (WrappedliquidstakedEther2.0)
$ws.Range("D15").Value = "1.894.61"
$ws.Range("E15").Value = "  -4.18%  "

# Row 16 (ShibaInu)
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008486"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.01%  "

# Row 17 (Litecoin)
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -13.34%  "

# Row 18 (WrappedBTC)
$ws.Range("D18").Value = "26.294.80"
$ws.Range("E18").Value = "  -5.90%  "

# Row 19 (Uniswap)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.920"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.34%  "

# Row 20 (Dai)
$ws.Range("E20").Value = "  +0.39%  "

# Row 21 (Avalanche)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.30%  "

# Row 22 (BitcoinCash)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "188.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.40%  "

# Row 23 (Chainlink)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.180"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.87%  "

# Row 24 (BinanceUSD)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.006"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.42%  "

# Row 25 (Monero)
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.54%  "

# Row 26 (Cosmos)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.664"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.63%  "

# Row 27 (Stellar)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1173"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.98%  "

# Row 28 (EthereumClassic)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.65"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.25%  "

# Row 29 (Hedera)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05817"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.25%  "

# Row 30 (Toncoin)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.279"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.77%  "

# Row 31 (PancakeSwap)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.324"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.95%  "

# Row 32 (InternetComputer(DFINITY))
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.529"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.65%  "

# Row 33 (Filecoin)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.507"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.38%  "

# Row 34 (LidoDAOToken)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.633"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.61%  "

# Row 35 (ARBITRUM)
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.010"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.81%  "

# Row 36 (ImmutableX)
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.5984"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.55%  "

# Row 37 (HuobiToken)
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.356"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.43%  "

# Row 38 (MXToken)
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.643"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.76%  "

# Row 39 (VeChain)
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01609"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.76%  "

# Row 40 (FraxShare)
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.013"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.36%  "

# Row 41 (Maker)
$ws.Range("D41").Value = "1.071.49"
$ws.Range("E41").Value = "  -4.66%  "

# Row 42 (TrustWalletToken)
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8628"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.40%  "

# Row 43 (PaxDollar)
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.009"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.47%  "

# Row 44 (Quant)
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.51%  "

# Row 45 (RocketPoolETH)
$ws.Range("D45").Value = "1.816.24"
$ws.Range("E45").Value = "  -3.90%  "

# Row 46 (BabyDogeCoin)
$ws.Range("E46").Value = "  +1.26%  "

# Row 47 (Aave)
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.10%  "

# Row 48 (Frax)
$ws.Range("E48").Value = "  +0.54%  "

# Row 49 (EnergySwap)
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.076"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.22%  "

# Row 50 (Mantle)
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4299"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.59%  "

# Row 51 (Cronos)
$ws.Range("E51").Value = "  -3.60%  "
